$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split the run "[onshow..now;frm=’yyyy-mm-dd hh:nn:ss’]" into five
#    runs, replacing the curly quotes (U+2019) with straight apostrophes.
# ---------------------------------------------------------------------

$range1 = $d.Content
$find = $range1.Find
$find.ClearFormatting()
$quote = [char]0x2019
$needle = "[onshow..now;frm=" + $quote + "yyyy-mm-dd hh:nn:ss" + $quote + "]"
$found = $find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $runStart = $range1.Start

    $piece1 = "[onshow..now;frm="
    $piece2 = "'"
    $piece3 = "yyyy-mm-dd hh:nn:ss"
    $piece4 = "'"
    $piece5 = "]"

    $s1 = $runStart
    $e1 = $s1 + $piece1.Length
    $e2 = $e1 + $piece2.Length
    $e3 = $e2 + $piece3.Length
    $e4 = $e3 + $piece4.Length
    $e5 = $e4 + $piece5.Length

    # Swap the curly quotes for straight ones in place (same length,
    # so no other offsets move).
    $q1 = $d.Range($e1, $e2)
    $q1.Text = $piece2
    $q2 = $d.Range($e3, $e4)
    $q2.Text = $piece4

    # Touch the formatting of each of the five segments so the engine
    # keeps them as separate runs instead of re-coalescing them (they
    # all share the same rPr, as in the target document).
    $segments = @(
        @($s1, $e1),
        @($e1, $e2),
        @($e2, $e3),
        @($e3, $e4),
        @($e4, $e5)
    )
    foreach ($seg in $segments) {
        $r = $d.Range($seg[0], $seg[1])
        $r.Bold = 1
        $r.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from right after the chart drawing to
#    right after the "[onshow..cst.PHP_VERSION]" run.
# ---------------------------------------------------------------------

$range2 = $d.Content
$find2 = $range2.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("[onshow..cst.PHP_VERSION]", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $target = $d.Range($range2.End, $range2.End)
    $d.Bookmarks.Add("_GoBack", $target)
}
